$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Copy formatting from the last existing data row (72) down to the new rows (73-84)
$ws.Range("A72:K72").Copy()
$ws.Range("A73:K84").PasteSpecial(-4122)

# Row 73: E0072
$ws.Range("A73").Value = 'E0072'
$ws.Range("B73").Value = 'Using Notice Layer and Icon in a Consent Notice'
$ws.Range("C73").Value = 'This example shows a Consent Notice that is structured in 2 layers - first for a summary and a second layer providing detailed overview of the intended processing. The layers contain controls for consent and rights which are accompanied by the label to be used e.g. for the button in the UI, and the icon to be displayed alongside it. This example also shows how a consent notice can be expressed in a machine-readable form in a manner that is similar to and can be used to create a graphical notice such as on a website or in a mobile app.'
$ws.Range("D73").Value = 'E0072.ttl'
$ws.Range("E73").Value = 'ttl'
$ws.Range("F73").Value = 'file'
$ws.Range("G73").Value = 'dpv:ConsentNotice,dpv:NoticeLayer,dpv:NoticeIcon,dpv:hasNoticeLayer,dpv:hasNoticeIcon'
$ws.Range("I73").Value = 'accepted'
$ws.Range("J73").Value = 45643
$ws.Range("K73").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(73).AutoFit()

# Row 74: E0073
$ws.Range("A74").Value = 'E0073'
$ws.Range("B74").Value = 'Tracking the status of a Notice across its lifecycle'
$ws.Range("C74").Value = 'In this example, statuses associated with a notice are maintained in a record. The first part shows a notice status within a common process representing all users showing the notice has been communicated i.e. shown or delivered to the users. The second method shows a data processing record where the notice status is included as part of other matters e.g. data collection or rights exercise. The third method shows a record containing only the notice status for a specific user. The fourth method shows how the notice status can also be recorded as part of the notice metadata itself.'
$ws.Range("D74").Value = 'E0073.ttl'
$ws.Range("E74").Value = 'ttl'
$ws.Range("F74").Value = 'file'
$ws.Range("G74").Value = 'dpv:NoticeStatus,dpv:hasNoticeStatus'
$ws.Range("I74").Value = 'accepted'
$ws.Range("J74").Value = 45643
$ws.Range("K74").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(74).AutoFit()

# Row 75: E0074
$ws.Range("A75").Value = 'E0074'
$ws.Range("B75").Value = 'Expressing involvement of a ''human subject'' in a process'
$ws.Range("C75").Value = 'This example shows how the involvement of humans can be expressed in various contexts. The first part shows how a specific individual can be asserted to be involved as a ''human subject'' in some process. The second part shows how a process or service can assert the involvement of specific categories of humans as subjects. The third part shows how humans can be expressed to be involved in a specific role (as participants) and are declared to be vulnerable. The fourth part clarifies the distinction between stating involvement of users as humans for oversight, and as data subjects in relation to their personal data.'
$ws.Range("D75").Value = 'E0074.ttl'
$ws.Range("E75").Value = 'ttl'
$ws.Range("F75").Value = 'file'
$ws.Range("G75").Value = 'dpv:HumanSubject,dpv:hasHumanSubject,dpv:DataSubject,dpv:hasDataSubject,dpv:hasHumanInvolvement,dpv:HumanInvolvement'
$ws.Range("I75").Value = 'accepted'
$ws.Range("J75").Value = 45643
$ws.Range("K75").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(75).AutoFit()

# Row 76: E0075
$ws.Range("A76").Value = 'E0075'
$ws.Range("B76").Value = 'Expressing involvement of tracking and profiling in processing'
$ws.Range("C76").Value = '<p>This example shows how the involvement of tracking and profiling can be expressed within a process. The first method simply states that tracking and profiling occur - this should be considered as bad practice as it does not clarify what the tracking and profiling involves in terms of data and entities - in particular what is being inferred as part of the profiling of data subjects. The second method is a verbose representation showing tracking is a composite process consisting of collecting and using location to generate a ''location tracking data'' for the individual, which is then used in the profiling process to infer financial status of the individual.</p>
<p>Through the explicit details of data involved in tracking and profiling, the transparency and accountability principles are easier to satisfy, and we can additionally perform appropriate risk/impact assessments to analyse whether the tracking and profiling should be allowed or not. The verbosity also shows that tracking and profiling are two separate processes - which is essential when the same tracking or profiling occurs in more than one process.  The last part shows how this verbosity can be reduced by taking the tracking and profiling (which as mentioned earlier are likely to be common to several processes) and assigning them a unique IRI or identifier, and using to indicate their involvement in a process.</p>'
$ws.Range("D76").Value = 'E0075.ttl'
$ws.Range("E76").Value = 'ttl'
$ws.Range("F76").Value = 'file'
$ws.Range("G76").Value = 'dpv:Tracking,dpv:Profiling'
$ws.Range("I76").Value = 'accepted'
$ws.Range("J76").Value = 45643
$ws.Range("K76").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(76).AutoFit()

# Row 77: E0076
$ws.Range("A77").Value = 'E0076'
$ws.Range("B77").Value = 'Representing contract metadata and controls'
$ws.Range("C77").Value = 'The below example shows using [[DCT]] and DPV properties to represent metadata about contracts. [[DCT]] is used for generic properties such as titles and descriptions, as well as contract-specific properties such as when it was accepted by all parties and the provenance of the contract document. The DPV properties are useful to express properties such as entities involved in specific roles such as Service Provider or a Data Controller, and to denote the type of contract which is useful for its interpretation - such as whether it is a negotiated or a standard form (non-negotiated) contract.'
$ws.Range("D77").Value = 'E0076.ttl'
$ws.Range("E77").Value = 'ttl'
$ws.Range("F77").Value = 'file'
$ws.Range("G77").Value = 'dpv:Contract'
$ws.Range("I77").Value = 'accepted'
$ws.Range("J77").Value = 45643
$ws.Range("K77").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(77).AutoFit()

# Row 78: E0077
$ws.Range("A78").Value = 'E0077'
$ws.Range("B78").Value = 'Representing the status of contracts'
$ws.Range("C78").Value = 'This example shows how the lifecycle of a contract in terms of its drafting and acceptance, as well as the fulfilment of its requirements and whether they have been breached.'
$ws.Range("D78").Value = 'E0077.ttl'
$ws.Range("E78").Value = 'ttl'
$ws.Range("F78").Value = 'file'
$ws.Range("G78").Value = 'dpv:Contract,dpv:ContractStatus,dpv:hasContractStatus,dpv:ContractFulfilmentStatus,dpv:hasContractualFulfilmentStatus'
$ws.Range("I78").Value = 'accepted'
$ws.Range("J78").Value = 45643
$ws.Range("K78").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(78).AutoFit()

# Row 79: E0078
$ws.Range("A79").Value = 'E0078'
$ws.Range("B79").Value = 'Representing clauses or terms within a contract'
$ws.Range("C79").Value = 'This example shows how specific clauses in a contract can be modelled, and also how their fulfilment status can be represented and tracked.'
$ws.Range("D79").Value = 'E0078.ttl'
$ws.Range("E79").Value = 'ttl'
$ws.Range("F79").Value = 'file'
$ws.Range("G79").Value = 'dpv:ContractualClause,dpv:hasContractualClause,dpv:ContractFulfilmentStatus,dpv:hasContractualFulfilmentStatus'
$ws.Range("I79").Value = 'accepted'
$ws.Range("J79").Value = 45643
$ws.Range("K79").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(79).AutoFit()

# Row 80: E0079
$ws.Range("A80").Value = 'E0079'
$ws.Range("B80").Value = 'Stating status of Legal Obligations'
$ws.Range("C80").Value = 'This example shows how the status of using legal obligations can be represented in terms of it being carried out and it being completed. Note that though the legal obligation has been stated as ''completed'', this is not an assessment of whether the requirements of the legal obligation have been completed and are verified, but instead the status represents the use of legal obligation as a legal basis within the process is completed.'
$ws.Range("D80").Value = 'E0079.ttl'
$ws.Range("E80").Value = 'ttl'
$ws.Range("F80").Value = 'file'
$ws.Range("G80").Value = 'dpv:LegalObligation,dpv:LegalObligationStatus'
$ws.Range("I80").Value = 'accepted'
$ws.Range("J80").Value = 45643
$ws.Range("K80").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(80).AutoFit()

# Row 81: E0080
$ws.Range("A81").Value = 'E0080'
$ws.Range("B81").Value = 'Stating status of Legitimate Interest'
$ws.Range("C81").Value = 'This example shows how the status of using legitimate interests can be represented in terms of whether the existence and use of legitimate interest has been communicated to the data subject, and whether the data subject has objected to its use.'
$ws.Range("D81").Value = 'E0080.ttl'
$ws.Range("E81").Value = 'ttl'
$ws.Range("F81").Value = 'file'
$ws.Range("G81").Value = 'dpv:LegitimateInterest,dpv:LegitimateInterestStatus'
$ws.Range("I81").Value = 'accepted'
$ws.Range("J81").Value = 45643
$ws.Range("K81").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(81).AutoFit()

# Row 82: E0081
$ws.Range("A82").Value = 'E0081'
$ws.Range("B82").Value = 'Stating status of using Official Authority'
$ws.Range("C82").Value = 'This example shows how the status of using the official authority as the legal basis can be represented.'
$ws.Range("D82").Value = 'E0081.ttl'
$ws.Range("E82").Value = 'ttl'
$ws.Range("F82").Value = 'file'
$ws.Range("G82").Value = 'dpv:OfficialAuthorityOfController,dpv:OfficialAuthorityExerciseStatus'
$ws.Range("I82").Value = 'accepted'
$ws.Range("J82").Value = 45643
$ws.Range("K82").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(82).AutoFit()

# Row 83: E0082
$ws.Range("A83").Value = 'E0082'
$ws.Range("B83").Value = 'Stating status of using Public Interest'
$ws.Range("C83").Value = 'This example shows how the status of using public interest as the legal basis can be represented.'
$ws.Range("D83").Value = 'E0082.ttl'
$ws.Range("E83").Value = 'ttl'
$ws.Range("F83").Value = 'file'
$ws.Range("G83").Value = 'dpv:PublicInterest,dpv:PublicInterestStatus'
$ws.Range("I83").Value = 'accepted'
$ws.Range("J83").Value = 45643
$ws.Range("K83").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(83).AutoFit()

# Row 84: E0083
$ws.Range("A84").Value = 'E0083'
$ws.Range("B84").Value = 'Stating status of using Vital Interest'
$ws.Range("C84").Value = 'This example shows how the status of using vital interest as the legal basis can be represented.'
$ws.Range("D84").Value = 'E0083.ttl'
$ws.Range("E84").Value = 'ttl'
$ws.Range("F84").Value = 'file'
$ws.Range("G84").Value = 'dpv:VitalInterest,dpv:VitalInterestStatus'
$ws.Range("I84").Value = 'accepted'
$ws.Range("J84").Value = 45643
$ws.Range("K84").Value = 'Harshvardhan J. Pandit'
$ws.Rows.Item(84).AutoFit()
